# Apply edits described by the commit "optimize dql test case logic" to the
# mysql_batchsql_cases.xlsx workbook (Sheet1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 17: optimize the query to use a schema-qualified table name ---
# (order matters for shared-string table layout: introduce new string
# values in the same order the original author's diff appended them)
$ws.Range("J17").Value = "select m.name as n from S.M as m order by m.name"
$ws.Range("G17").Value = "S.M"
$ws.Range("F17").Value = "S"
$ws.Range("E17").Value = "Schema"

# --- Rows 2-16: Testable column (B) flips from "y" to "n" ---
foreach ($r in 2..16) {
    $ws.Cells.Item($r, 2).Value = "n"
}

# --- Selection moves to B22 ---
$ws.Range("B22").Select()
